# Applies the "Ready for final check" commit to Eco_GG_soft.xlsx
#
# Summary of the data edits (everything else in the diff is either
# shared-string renumbering caused by these edits, or cosmetic
# Excel-version metadata that carries no semantic content):
#
#   Sheet "kite":
#     - Rows 1-7 get new parameter names / values (structure.fixed.* ->
#       structure.soft.* plus a couple of renamed onboard-systems rows).
#     - Rows 8-12 (now redundant / superseded) are removed entirely.
#
#   Sheet "gStation":
#     - hydMotor.p_1 (B29): 0 -> 200
#     - hydMotor.p_2 (B30): 0 -> 80

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "kite"
# ---------------------------------------------------------------
$kite = $wb.Worksheets.Item("kite")

# Drop the old rows 8:12 completely (content + formatting) so no
# leftover <row> elements remain for them.
$kite.Range("A8:B12").Clear()

# Row 1: structure.fixed.approach (1) -> structure.soft.p_fabric (45)
$kite.Range("A1").Value = "structure.soft.p_fabric"
$kite.Range("B1").Value = 45

# Row 2: structure.fixed.one.p_str (250) -> structure.soft.p_bridle (8)
$kite.Range("A2").Value = "structure.soft.p_bridle"
$kite.Range("B2").Value = 8

# Row 3: structure.fixed.one.p_wet (200) -> structure.soft.L_str (0.57)
$kite.Range("A3").Value = "structure.soft.L_str"
$kite.Range("B3").Value = 0.57

# Row 4: structure.fixed.two.p_uni (3) -> obgen.p (120)
$kite.Range("A4").Value = "obgen.p"
$kite.Range("B4").Value = 120

# Row 5: structure.fixed.two.p_tri (3.6) -> prop.p (120)
$kite.Range("A5").Value = "prop.p"
$kite.Range("B5").Value = 120

# Row 6: structure.fixed.two.f_man (0.75) -> obBatt.p (150) -- keeps its
# existing number format.
$kite.Range("A6").Value = "obBatt.p"
$kite.Range("B6").Value = 150

# Row 7: structure.soft.p_A (45) -> avio.C (150000), now carrying the
# same number format as the old row 12 (avionics.C) that it replaces.
$kite.Range("A7").Value = "avio.C"
$kite.Range("B7").Value = 150000
$kite.Range("B7").NumberFormat = $kite.Range("B6").NumberFormat

# ---------------------------------------------------------------
# Sheet "gStation"
# ---------------------------------------------------------------
$gStation = $wb.Worksheets.Item("gStation")

$gStation.Range("B29").Value = 200
$gStation.Range("B30").Value = 80
